$d = $word.ActiveDocument

# --- 1 & 2: In the "Bar01.mat, RecsAnalyze02_ba1 ..." bullet, the runs were
# re-typed/re-saved by the author without any textual change, which causes
# Word to recombine adjacent runs that share identical formatting. We
# reproduce that by doing an in-place (identity) Find/Replace on the first
# run's text; Word then renormalizes run boundaries for the rest of the
# paragraph as a side effect, merging the bold "Bar01.mat," / " RecsAnalyze02_ba1 "
# runs into one, and the three non-bold runs that follow into another.
$d.Content.Find.Execute("Bar01.mat,", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Bar01.mat,", 2) | Out-Null

# --- 3: The header's cached PAGE field result changes from "2" to "1".
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1", 2) | Out-Null
